$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 260
$ws.Range("I15").Value = 260
$ws.Range("K15").Value = 780
$ws.Range("M15").Value = -611

$ws.Range("H28").Value = 451.94446
$ws.Range("J28").Value = 223.75
$ws.Range("L28").Value = 223.75
$ws.Range("N28").Value = -1193.75

$ws.Range("H98").Value = 1024.7826
$ws.Range("I98").Value = 1027.1428
$ws.Range("K98").Value = 1027.1428
$ws.Range("M98").Value = 470.8571999999999

$ws.Range("H107").Value = 511.6842
$ws.Range("I107").Value = 534.9231
$ws.Range("J107").Value = 461.33334
$ws.Range("K107").Value = 534.9231
$ws.Range("L107").Value = 461.33334
$ws.Range("M107").Value = 1385.0769
$ws.Range("N107").Value = -4301.33334

$ws.Range("H122").Value = 1024.7826
$ws.Range("I122").Value = 1027.1428
$ws.Range("K122").Value = 3081.4284
$ws.Range("M122").Value = -631.4284000000002

$ws.Range("H125").Value = 1067.75
$ws.Range("I125").Value = 1173.8334
$ws.Range("K125").Value = 10564.5006
$ws.Range("M125").Value = -8104.500599999999

$ws.Range("H129").Value = 699709.25
$ws.Range("J129").Value = 862362.0600000001
$ws.Range("L129").Value = 2587086.18
$ws.Range("N129").Value = -2597086.18

$ws.Range("H132").Value = 1604.2439
$ws.Range("I132").Value = 1631.85
$ws.Range("J132").Value = 500
$ws.Range("K132").Value = 4895.549999999999
$ws.Range("L132").Value = 1500
$ws.Range("M132").Value = -2365.549999999999
$ws.Range("N132").Value = -6560

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1415.1082
$ws.Range("I61").Value = 1418.3055
$ws.Range("J61").Value = 1300
$ws.Range("K61").Value = 1418.3055
$ws.Range("L61").Value = 1300
$ws.Range("M61").Value = -1206.3055
$ws.Range("N61").Value = -1724

$ws.Range("H74").Value = 1020.2917
$ws.Range("I74").Value = 1016.8261
$ws.Range("K74").Value = 1016.8261
$ws.Range("M74").Value = -142.8261

$ws.Range("H77").Value = 1020.2917
$ws.Range("I77").Value = 1016.8261
$ws.Range("K77").Value = 5084.1305
$ws.Range("M77").Value = -716.1305000000002

$ws.Range("H132").Value = 1970.2444
$ws.Range("I132").Value = 1676.225
$ws.Range("K132").Value = 5028.674999999999
$ws.Range("M132").Value = -2498.674999999999

$ws.Range("H136").Value = 1415.1082
$ws.Range("I136").Value = 1418.3055
$ws.Range("J136").Value = 1300
$ws.Range("K136").Value = 4254.916499999999
$ws.Range("L136").Value = 3900
$ws.Range("M136").Value = -1704.916499999999
$ws.Range("N136").Value = -9000

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 19677.701
$ws.Range("I134").Value = 1800.5641
$ws.Range("K134").Value = 5401.692300000001
$ws.Range("M134").Value = -2866.692300000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1257.3889
$ws.Range("I16").Value = 1168.5
$ws.Range("J16").Value = 1435.1666
$ws.Range("K16").Value = 1168.5
$ws.Range("L16").Value = 1435.1666
$ws.Range("M16").Value = -881.5
$ws.Range("N16").Value = -2009.1666

$ws.Range("H58").Value = 2864.386
$ws.Range("I58").Value = 1057.1794
$ws.Range("J58").Value = 6780
$ws.Range("K58").Value = 1057.1794
$ws.Range("L58").Value = 6780
$ws.Range("M58").Value = -854.1794
$ws.Range("N58").Value = -7186

$ws.Range("H94").Value = 6317.4287
$ws.Range("J94").Value = 6883.2
$ws.Range("L94").Value = 6883.2
$ws.Range("N94").Value = -7785.2

$ws.Range("H100").Value = 25390
$ws.Range("J100").Value = 25390
$ws.Range("L100").Value = 25390
$ws.Range("N100").Value = -27554

$ws.Range("H113").Value = 1257.3889
$ws.Range("I113").Value = 1168.5
$ws.Range("J113").Value = 1435.1666
$ws.Range("K113").Value = 1168.5
$ws.Range("L113").Value = 1435.1666
$ws.Range("M113").Value = 1001.5
$ws.Range("N113").Value = -5775.1666

$ws.Range("H132").Value = 1691
$ws.Range("I132").Value = 1032.4117
$ws.Range("J132").Value = 3090.5
$ws.Range("K132").Value = 3097.2351
$ws.Range("L132").Value = 9271.5
$ws.Range("M132").Value = -567.2351000000003
$ws.Range("N132").Value = -14331.5

$ws.Range("H134").Value = 15626167
$ws.Range("I134").Value = 1187.3684
$ws.Range("K134").Value = 3562.1052
$ws.Range("M134").Value = -1027.1052

$ws.Range("H136").Value = 2864.386
$ws.Range("I136").Value = 1057.1794
$ws.Range("J136").Value = 6780
$ws.Range("K136").Value = 3171.5382
$ws.Range("L136").Value = 20340
$ws.Range("M136").Value = -621.5382
$ws.Range("N136").Value = -25440

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 672.6923
$ws.Range("I5").Value = 316
$ws.Range("J5").Value = 1088.8334
$ws.Range("K5").Value = 948
$ws.Range("L5").Value = 3266.5002
$ws.Range("M5").Value = -836
$ws.Range("N5").Value = -3490.5002

$ws.Range("H107").Value = 389707.34
$ws.Range("J107").Value = 556285.5
$ws.Range("L107").Value = 1668856.5
$ws.Range("N107").Value = -1672696.5

$ws.Range("H135").Value = 672.6923
$ws.Range("I135").Value = 316
$ws.Range("J135").Value = 1088.8334
$ws.Range("K135").Value = 2844
$ws.Range("L135").Value = 9799.500599999999
$ws.Range("M135").Value = -309
$ws.Range("N135").Value = -14869.5006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1828.3529
$ws.Range("I113").Value = 1007.46155
$ws.Range("J113").Value = 4496.25
$ws.Range("K113").Value = 1007.46155
$ws.Range("L113").Value = 4496.25
$ws.Range("M113").Value = 1162.53845
$ws.Range("N113").Value = -8836.25

$ws.Range("H126").Value = 2574.077
$ws.Range("I126").Value = 2766.3
$ws.Range("J126").Value = 1933.3334
$ws.Range("K126").Value = 8298.900000000001
$ws.Range("L126").Value = 5800.0002
$ws.Range("M126").Value = -5828.900000000001
$ws.Range("N126").Value = -10740.0002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2717.7727
$ws.Range("I132").Value = 2536.75
$ws.Range("J132").Value = 3200.5
$ws.Range("K132").Value = 7610.25
$ws.Range("L132").Value = 9601.5
$ws.Range("M132").Value = -5080.25
$ws.Range("N132").Value = -14661.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1281.1904
$ws.Range("J132").Value = 1539.1428
$ws.Range("L132").Value = 4617.428400000001
$ws.Range("N132").Value = -9677.428400000001

$ws.Range("H136").Value = 856.14703
$ws.Range("I136").Value = 646.8889
$ws.Range("J136").Value = 1663.2858
$ws.Range("K136").Value = 1940.6667
$ws.Range("L136").Value = 4989.857400000001
$ws.Range("M136").Value = 609.3332999999998
$ws.Range("N136").Value = -10089.8574
